$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 9 and 8 (bottom-up so row indices of earlier rows are unaffected)
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()

# Row 2: FAPs | Wnt5a | Ror1 | ECs  (updated M:T values)
$ws.Cells.Item(2,13).Value = 0.1350193333333333
$ws.Cells.Item(2,14).Value = 0.405058
$ws.Cells.Item(2,15).Value = 0.004770957739100914
$ws.Cells.Item(2,16).Value = 0.004770957739100915
$ws.Cells.Item(2,17).Value = 1.408824578704444
$ws.Cells.Item(2,18).Value = 12.67942120834
$ws.Cells.Item(2,19).Value = 0.004633120943133651
$ws.Cells.Item(2,20).Value = 0.004633120943133653

# Row 3: FAPs | Wnt5a | Ror1 | FAPs (updated O,P,S,T values)
$ws.Cells.Item(3,15).Value = 0.4358572565646723
$ws.Cells.Item(3,16).Value = 0.4358572565646723
$ws.Cells.Item(3,19).Value = 0.4232649908123294
$ws.Cells.Item(3,20).Value = 0.4232649908123295

# Row 4: FAPs | Wnt5a | Ror1 | MuSCs (updated M:T values)
$ws.Cells.Item(4,13).Value = 15.83036566666667
$ws.Cells.Item(4,14).Value = 47.491097
$ws.Cells.Item(4,15).Value = 0.5593717856962268
$ws.Cells.Item(4,16).Value = 0.5593717856962268
$ws.Cells.Item(4,17).Value = 165.1778874216455
$ws.Cells.Item(4,18).Value = 1486.60098679481
$ws.Cells.Item(4,19).Value = 0.5432110861236953
$ws.Cells.Item(4,20).Value = 0.5432110861236953

# Row 5: was FAPs|Wnt5a|Ror1|Resolving-Mac -> now MuSCs|Wnt5a|Ror1|ECs
$ws.Cells.Item(5,1).Value = "MuSCs"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.310422
$ws.Cells.Item(5,8).Value = 0.9312659999999999
$ws.Cells.Item(5,9).Value = 0.02889080212084161
$ws.Cells.Item(5,10).Value = 0.02889080212084161
$ws.Cells.Item(5,13).Value = 0.1350193333333333
$ws.Cells.Item(5,14).Value = 0.405058
$ws.Cells.Item(5,15).Value = 0.004770957739100914
$ws.Cells.Item(5,16).Value = 0.004770957739100915
$ws.Cells.Item(5,17).Value = 0.041912971492
$ws.Cells.Item(5,18).Value = 0.3772167434279999
$ws.Cells.Item(5,19).Value = 0.0001378367959672624
$ws.Cells.Item(5,20).Value = 0.0001378367959672624

# Row 6: was MuSCs|Wnt5a|Ror1|ECs -> now MuSCs|Wnt5a|Ror1|FAPs
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 12.33487266666667
$ws.Cells.Item(6,14).Value = 37.004618
$ws.Cells.Item(6,15).Value = 0.4358572565646723
$ws.Cells.Item(6,16).Value = 0.4358572565646723
$ws.Cells.Item(6,17).Value = 3.829015842932
$ws.Cells.Item(6,18).Value = 34.461142586388
$ws.Cells.Item(6,19).Value = 0.01259226575234284
$ws.Cells.Item(6,20).Value = 0.01259226575234284

# Row 7: was MuSCs|Wnt5a|Ror1|FAPs -> now MuSCs|Wnt5a|Ror1|MuSCs
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,13).Value = 15.83036566666667
$ws.Cells.Item(7,14).Value = 47.491097
$ws.Cells.Item(7,15).Value = 0.5593717856962268
$ws.Cells.Item(7,16).Value = 0.5593717856962268
$ws.Cells.Item(7,17).Value = 4.914093770977999
$ws.Cells.Item(7,18).Value = 44.22684393880199
$ws.Cells.Item(7,19).Value = 0.01616069957253151
$ws.Cells.Item(7,20).Value = 0.01616069957253151

$wb.Save()
